# Add columns I0 (I) and IF (J) to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy formatting (bold, border, center) from H1, then set text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for rows 2..75: column I (I0) and column J (IF)
$data = @(
    @(2, 2),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(5, 5),
    @(7, 7),
    @(6, 6),
    @(8, 8),
    @(7, 7),
    @(8, 9),
    @(6, 6),
    @(5, 6),
    @(2, 2),
    @(4, 4),
    @(7, 7),
    @(1, 2),
    @(8, 8),
    @(5, 5),
    @(7, 7),
    @(7, 8),
    @(2, 3),
    @(6, 6),
    @(6, 7),
    @(6, 6),
    @(6, 6),
    @(5, 5),
    @(8, 8),
    @(9, 9),
    @(5, 6),
    @(5, 5),
    @(3, 4),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(4, 4),
    @(8, 8),
    @(5, 6),
    @(6, 6),
    @(6, 6),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(7, 8),
    @(6, 6),
    @(5, 6),
    @(7, 7),
    @(6, 6),
    @(7, 8),
    @(8, 8),
    @(9, 9),
    @(8, 9),
    @(8, 8),
    @(8, 8),
    @(5, 6),
    @(8, 8),
    @(5, 6),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(10, 10),
    @(8, 8),
    @(5, 5),
    @(4, 4),
    @(5, 6),
    @(8, 8),
    @(7, 7)
)

for ($r = 2; $r -le 75; $r++) {
    $pair = $data[$r - 2]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}
